# Scrum Master: Files update
# Remove the completed task "Começar a fazer o use case diagram" from the
# "week1" sheet (first sheet, B12) and move the selection to that cell,
# matching the user having selected B12 and pressed Delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")
$ws.Activate()

# Clear the contents of B12 (removes the shared-string reference, the row
# collapses out of sheetData on save since it becomes empty, but later rows
# are not shifted).
$ws.Range("B12").ClearContents()

# Leave the selection on the cell that was just cleared.
$ws.Range("B12").Select()
